$wb = $excel.ActiveWorkbook

# --- Sheet ALC (hunk 0) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1682.4445
$ws.Range("I15").Value = 1682.4445
$ws.Range("K15").Value = 5047.333500000001
$ws.Range("M15").Value = -4878.333500000001

# --- Sheet ALC (hunk 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 7323
$ws.Range("J17").Value = 7323
$ws.Range("L17").Value = 21969
$ws.Range("N17").Value = -22305

# --- Sheet ALC (hunk 2) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 30303754
$ws.Range("I28").Value = 40000520
$ws.Range("K28").Value = 40000520
$ws.Range("M28").Value = -40000035

# --- Sheet ALC (hunk 3) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 16391.75
$ws.Range("J64").Value = 9100
$ws.Range("L64").Value = 9100
$ws.Range("N64").Value = -9596

# --- Sheet ALC (hunk 4) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 16391.75
$ws.Range("J67").Value = 9100
$ws.Range("L67").Value = 9100
$ws.Range("N67").Value = -10816

# --- Sheet ALC (hunk 5) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1363.9131
$ws.Range("I98").Value = 1497.8422
$ws.Range("J98").Value = 727.75
$ws.Range("K98").Value = 1497.8422
$ws.Range("L98").Value = 727.75
$ws.Range("M98").Value = 0.157799999999952
$ws.Range("N98").Value = -3723.75

# --- Sheet ALC (hunk 6) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1921.4286
$ws.Range("J112").Value = 1962.963
$ws.Range("L112").Value = 5888.889
$ws.Range("N112").Value = -8104.889

# --- Sheet ALC (hunk 7) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1363.9131
$ws.Range("I122").Value = 1497.8422
$ws.Range("J122").Value = 727.75
$ws.Range("K122").Value = 4493.5266
$ws.Range("L122").Value = 2183.25
$ws.Range("M122").Value = -2043.5266
$ws.Range("N122").Value = -7083.25

# --- Sheet ALC (hunk 8) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 12130.962
$ws.Range("I132").Value = 5824.1304
$ws.Range("K132").Value = 17472.3912
$ws.Range("M132").Value = -14942.3912

# --- Sheet ALC (hunk 9) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2755.6038
$ws.Range("I137").Value = 2564.625
$ws.Range("K137").Value = 7693.875
$ws.Range("M137").Value = -5143.875

# --- Sheet ALC (hunk 10) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3732.525
$ws.Range("I138").Value = 1227.7693
$ws.Range("K138").Value = 3683.3079
$ws.Range("M138").Value = 1456.6921

# --- Sheet ALC (hunk 11) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 91999.8
$ws.Range("I139").Value = 59999
$ws.Range("K139").Value = 59999
$ws.Range("M139").Value = -54859

# --- Sheet ALC (hunk 12) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3887.8
$ws.Range("I141").Value = 2134.6453
$ws.Range("J141").Value = 17474.75
$ws.Range("K141").Value = 6403.9359
$ws.Range("L141").Value = 52424.25
$ws.Range("M141").Value = -1223.9359
$ws.Range("N141").Value = -62784.25

# --- Sheet ARM (hunk 13) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 391.66666
$ws.Range("J3").Value = 375
$ws.Range("L3").Value = 375
$ws.Range("N3").Value = -605

# --- Sheet ARM (hunk 14) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5731.353
$ws.Range("I32").Value = 5862.625
$ws.Range("K32").Value = 5862.625
$ws.Range("M32").Value = -5575.625

# --- Sheet ARM (hunk 15) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2585.6
$ws.Range("I61").Value = 1116.8182
$ws.Range("K61").Value = 1116.8182
$ws.Range("M61").Value = -904.8181999999999

# --- Sheet ARM (hunk 16) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1480.6666
$ws.Range("I74").Value = 1214.7858
$ws.Range("K74").Value = 1214.7858
$ws.Range("M74").Value = -340.7858000000001

# --- Sheet ARM (hunk 17) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1480.6666
$ws.Range("I77").Value = 1214.7858
$ws.Range("K77").Value = 6073.929
$ws.Range("M77").Value = -1705.929

# --- Sheet ARM (hunk 18) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1266.6364
$ws.Range("I132").Value = 1226.5
$ws.Range("K132").Value = 3679.5
$ws.Range("M132").Value = -1149.5

# --- Sheet ARM (hunk 19) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2585.6
$ws.Range("I136").Value = 1116.8182
$ws.Range("K136").Value = 3350.4546
$ws.Range("M136").Value = -800.4546

# --- Sheet BSM (hunk 20) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4806.9375
$ws.Range("I86").Value = 1800
$ws.Range("J86").Value = 5007.4
$ws.Range("K86").Value = 1800
$ws.Range("L86").Value = 5007.4
$ws.Range("M86").Value = -677
$ws.Range("N86").Value = -7253.4

# --- Sheet BSM (hunk 21) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4806.9375
$ws.Range("I89").Value = 1800
$ws.Range("J89").Value = 5007.4
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 25037
$ws.Range("M89").Value = -3384
$ws.Range("N89").Value = -36269

# --- Sheet CRP (hunk 22) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 1189.4
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1189.4
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 1189.4
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -1469.4

# --- Sheet CRP (hunk 23) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1929.8334
$ws.Range("I31").Value = 1779.1333
$ws.Range("J31").Value = 2683.3333
$ws.Range("K31").Value = 1779.1333
$ws.Range("L31").Value = 2683.3333
$ws.Range("M31").Value = -1484.1333
$ws.Range("N31").Value = -3273.3333

# --- Sheet CRP (hunk 24) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1929.8334
$ws.Range("I34").Value = 1779.1333
$ws.Range("J34").Value = 2683.3333
$ws.Range("K34").Value = 1779.1333
$ws.Range("L34").Value = 2683.3333
$ws.Range("M34").Value = -1577.1333
$ws.Range("N34").Value = -3087.3333

# --- Sheet CRP (hunk 25) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2684.25
$ws.Range("I58").Value = 2030.1428
$ws.Range("J58").Value = 3600
$ws.Range("K58").Value = 2030.1428
$ws.Range("L58").Value = 3600
$ws.Range("M58").Value = -1827.1428
$ws.Range("N58").Value = -4006

# --- Sheet CRP (hunk 26) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 13128.125
$ws.Range("I103").Value = 13128.125
$ws.Range("K103").Value = 13128.125
$ws.Range("M103").Value = -11956.125

# --- Sheet CRP (hunk 27) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2049.0942
$ws.Range("I132").Value = 1911.0667
$ws.Range("K132").Value = 5733.2001
$ws.Range("M132").Value = -3203.2001

# --- Sheet CRP (hunk 28) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3980.9355
$ws.Range("I134").Value = 3181.524
$ws.Range("K134").Value = 9544.572
$ws.Range("M134").Value = -7009.572

# --- Sheet CRP (hunk 29) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2684.25
$ws.Range("I136").Value = 2030.1428
$ws.Range("J136").Value = 3600
$ws.Range("K136").Value = 6090.428400000001
$ws.Range("L136").Value = 10800
$ws.Range("M136").Value = -3540.428400000001
$ws.Range("N136").Value = -15900

# --- Sheet CUL (hunk 30) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 67.85714
$ws.Range("I14").Value = 67.85714
$ws.Range("K14").Value = 203.57142
$ws.Range("M14").Value = -30.57141999999999

# --- Sheet CUL (hunk 31) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").ClearContents()

# --- Sheet CUL (hunk 32) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 2051.2354
$ws.Range("I103").Value = 1490.0834
$ws.Range("J103").Value = 3398
$ws.Range("K103").Value = 4470.2502
$ws.Range("L103").Value = 10194
$ws.Range("M103").Value = -3591.2502
$ws.Range("N103").Value = -11952

# --- Sheet CUL (hunk 33) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 4770.875
$ws.Range("I118").Value = 1084
$ws.Range("K118").Value = 3252
$ws.Range("M118").Value = -2009

# --- Sheet CUL (hunk 34) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3310.818
$ws.Range("J131").Value = 4370.154
$ws.Range("L131").Value = 13110.462
$ws.Range("N131").Value = -23190.462

# --- Sheet CUL (hunk 35) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 3281.6667
$ws.Range("I133").Value = 3281.6667
$ws.Range("K133").Value = 9845.000100000001
$ws.Range("M133").Value = -4785.000100000001

# --- Sheet GSM (hunk 36) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 12389
$ws.Range("I41").Value = 1500
$ws.Range("K41").Value = 1500
$ws.Range("M41").Value = -1145

# --- Sheet GSM (hunk 37) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 33447536
$ws.Range("I80").Value = 558055.5
$ws.Range("J80").Value = 41669904
$ws.Range("K80").Value = 558055.5
$ws.Range("L80").Value = 41669904
$ws.Range("M80").Value = -557057.5
$ws.Range("N80").Value = -41671900

# --- Sheet GSM (hunk 38) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 33447536
$ws.Range("I83").Value = 558055.5
$ws.Range("J83").Value = 41669904
$ws.Range("K83").Value = 2790277.5
$ws.Range("L83").Value = 208349520
$ws.Range("M83").Value = -2785285.5
$ws.Range("N83").Value = -208359504

# --- Sheet GSM (hunk 39) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2440.524
$ws.Range("I132").Value = 2302.4722
$ws.Range("K132").Value = 6907.4166
$ws.Range("M132").Value = -4377.4166

# --- Sheet LTW (hunk 40) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4845.0625
$ws.Range("I7").Value = 4181.6
$ws.Range("K7").Value = 4181.6
$ws.Range("M7").Value = -4069.6

# --- Sheet LTW (hunk 41) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# --- Sheet LTW (hunk 42) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4845.0625
$ws.Range("I126").Value = 4181.6
$ws.Range("K126").Value = 12544.8
$ws.Range("M126").Value = -10074.8

# --- Sheet LTW (hunk 43) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4137.6787
$ws.Range("I132").Value = 2004.5
$ws.Range("K132").Value = 6013.5
$ws.Range("M132").Value = -3483.5

# --- Sheet LTW (hunk 44) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7778.3335
$ws.Range("I136").Value = 7799.8
$ws.Range("J136").Value = 7751.5
$ws.Range("K136").Value = 23399.4
$ws.Range("L136").Value = 23254.5
$ws.Range("M136").Value = -20849.4
$ws.Range("N136").Value = -28354.5

# --- Sheet WVR (hunk 45) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 10000
$ws.Range("J42").Value = 10000
$ws.Range("L42").Value = 10000
$ws.Range("N42").Value = -10756

# --- Sheet WVR (hunk 46) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2424.2144
$ws.Range("I126").Value = 2528.3333
$ws.Range("J126").Value = 1799.5
$ws.Range("K126").Value = 7584.999899999999
$ws.Range("L126").Value = 5398.5
$ws.Range("M126").Value = -5114.999899999999
$ws.Range("N126").Value = -10338.5

# --- Sheet WVR (hunk 47) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1947.5
$ws.Range("I136").Value = 1056.5
$ws.Range("K136").Value = 3169.5
$ws.Range("M136").Value = -619.5
